$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.66000000000057
$ws.Range("H2").Value = 0.0000000000000005621382403165349
$ws.Range("K2").Value = 43.49676309937972
$ws.Range("L2").Value = "[35.40663208792838, 51.58689411083106]"
$ws.Range("O2").Value = 1.478026573760964
$ws.Range("P2").Value = "[1.2767633807381937, 1.679289766783734]"
$ws.Range("S2").Value = 62.64349526615166
$ws.Range("T2").Value = "[57.44679274818219, 67.84019778412113]"
$ws.Range("W2").Value = 19.6238638638643
$ws.Range("X2").Value = 18.80192192192234
$ws.Range("Y2").Value = 20.44580580580626

# Row 3 updates
$ws.Range("E3").Value = 23.34000000000021
$ws.Range("G3").Value = 0.000000006638526839353176
$ws.Range("H3").Value = 0.00000002128374678146951
$ws.Range("K3").Value = 36.52033460028461
$ws.Range("L3").Value = "[22.54488015608607, 50.495789044483146]"
$ws.Range("M3").Value = 0.0000008004960514451653
$ws.Range("N3").Value = 0.0000008004960514451653
$ws.Range("O3").Value = -2.679316257115619
$ws.Range("P3").Value = "[-3.094421592725082, -2.264210921506157]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 61.69649872702855
$ws.Range("T3").Value = "[54.01057817449286, 69.38241927956423]"
$ws.Range("W3").Value = 9.952792792792883
$ws.Range("X3").Value = 8.410810810810885
$ws.Range("Y3").Value = 11.49477477477488
